$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Paragraph 1: "There should be two hyperlinks below:"
#    -> "There should be two hyperlinks below, the first one in red:"
#    and move the "_GoBack" bookmark from the last (now-empty) paragraph
#    to sit right before the final ":" run.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("below")
if ($found) {
    $rng.Collapse(0)                      # wdCollapseEnd
    $bm = $d.Bookmarks("_GoBack")
    $bm.Delete()
    $d.Bookmarks.Add("_GoBack", $rng)
    $rng.InsertAfter(", the first one in red")
}

# ------------------------------------------------------------------
# 2) First "+++LINK" paragraph (the one with a label) gets:
#    - "+++LINK " and "({ " runs merged into a single run
#    - the gramStart/gramEnd proofErr markers removed (spellStart/spellEnd kept)
#    - every run (and the paragraph mark) coloured red (FF0000)
# ------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$start = $p2.Range.Start
$end = $p2.Range.End
$r2 = $d.Range($start, $end)
$xml2 = '<w:p>' +
    '<w:r><w:t xml:space="preserve">+++LINK ({ </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>url</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    "<w:r><w:t>: 'http://www.apple.com', label: 'Apple' })+++</w:t></w:r>" +
    '</w:p>'
$r2.InsertXML($xml2)

$p2b = $d.Paragraphs.Item(2)
$p2b.Range.Font.Color = 255               # wdColorRed -> <w:color w:val="FF0000"/>

# ------------------------------------------------------------------
# 3) Second "+++LINK" paragraph (no label) gets:
#    - "+++LINK " and "({ " runs merged into a single run
#    - the gramStart/gramEnd proofErr markers removed (spellStart/spellEnd kept)
#    - no colour change
# ------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$start3 = $p3.Range.Start
$end3 = $p3.Range.End
$r3 = $d.Range($start3, $end3)
$xml3 = '<w:p>' +
    '<w:r><w:t xml:space="preserve">+++LINK ({ </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>url</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    "<w:r><w:t>: 'http://www.apple.com' })+++</w:t></w:r>" +
    '</w:p>'
$r3.InsertXML($xml3)
